$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 148, shifting existing rows 148:171 down to 149:172
$ws.Rows.Item(148).Insert()

# Populate the new row 148 with the new data record
$ws.Cells.Item(148, 1).Value = 11
$ws.Cells.Item(148, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(148, 3).Value = "Bíobío"
$ws.Cells.Item(148, 4).Value = 44504
$ws.Cells.Item(148, 5).Value = 8
$ws.Cells.Item(148, 6).Value = 100112023
$ws.Cells.Item(148, 7).Value = "Brócoli"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 2200
$ws.Cells.Item(148, 11).Value = 650
$ws.Cells.Item(148, 12).Value = 700
$ws.Cells.Item(148, 13).Value = 673
$ws.Cells.Item(148, 14).Value = "`$/unidad"
$ws.Cells.Item(148, 15).Value = "Región Metropolitana"
$ws.Cells.Item(148, 16).Value = 673
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = "Hortaliza"
